# Card2: fix "Correction " header -> "Correction", add new "Serviced by "
# column (O) with matching header style, and extend the N/O data columns
# so every data row (2..13) carries an explicit (empty) text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# --- Header row ---------------------------------------------------------
# N1 had a trailing space ("Correction ") -> drop it.
$ws.Cells.Item(1, 14).Value = "Correction"

# New header O1, styled like the other header cells (bold/border/center).
$ws.Cells.Item(1, 15).Value = "Serviced by "
$ws.Cells.Item(1, 14).Copy()
$ws.Cells.Item(1, 15).PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows -----------------------------------------------------------
for ($r = 2; $r -le 13; $r++) {
    # N2:N13 were present but empty -> they now hold the literal text "nan"
    # (matching every other "nan" placeholder cell on this sheet).
    $ws.Cells.Item($r, 14).Value = "nan"

    # O2:O13 are new cells that stay blank, but need to be real empty TEXT
    # cells (not simply absent). Writing "'" forces a text cell with empty
    # content, then resetting the style keeps it on the default style (no
    # quote-prefix formatting carried over).
    $ws.Cells.Item($r, 15).Value = "'"
    $ws.Cells.Item($r, 15).Style = "Normal"
}
